# Update "想去人数" (want-to-go count) figures in the 展览 (sheet 1) and
# 全部类型 (sheet 4) worksheets to match the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsExhibit.Range("F3").Value = 1451
$wsExhibit.Range("F5").Value = 12108
$wsExhibit.Range("F6").Value = 4472
$wsExhibit.Range("F7").Value = 41
$wsExhibit.Range("F13").Value = 182
$wsExhibit.Range("F15").Value = 5241
$wsExhibit.Range("F16").Value = 64
$wsExhibit.Range("F18").Value = 545
$wsExhibit.Range("F20").Value = 11481

$wsAll = $wb.Worksheets.Item(4)       # 全部类型
$wsAll.Range("F3").Value = 1451
$wsAll.Range("F5").Value = 12108
$wsAll.Range("F6").Value = 4472
$wsAll.Range("F7").Value = 41
$wsAll.Range("F14").Value = 182
$wsAll.Range("F16").Value = 5241
$wsAll.Range("F17").Value = 64
$wsAll.Range("F19").Value = 545
$wsAll.Range("F21").Value = 11481
